$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I3").Value = -0.3481849190473605
$ws.Range("J3").Value = 0.2123102164198775
$ws.Range("K3").Value = -0.6323858570326343
$ws.Range("L3").Value = 2.935572397170822
